$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-09-02 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-03 Wednesday", 2)

# Update each division problem cell by its (row, column) position in the table.
# Each replacement is scoped to the individual cell's Range and uses wdReplaceOne (1)
# so that duplicate expressions (e.g. "62div5=") occurring in more than one cell are
# each replaced independently with their own target value, instead of a document-wide
# wdReplaceAll clobbering every matching occurrence.
$tbl = $d.Tables.Item(1)

$cell = $tbl.Cell(1, 1)
$cell.Range.Find.Execute("88÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "90÷8=", 1)

$cell = $tbl.Cell(1, 2)
$cell.Range.Find.Execute("62÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷9=", 1)

$cell = $tbl.Cell(1, 3)
$cell.Range.Find.Execute("47÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "24÷7=", 1)

$cell = $tbl.Cell(1, 4)
$cell.Range.Find.Execute("64÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "91÷8=", 1)

$cell = $tbl.Cell(1, 5)
$cell.Range.Find.Execute("85÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "36÷9=", 1)

$cell = $tbl.Cell(5, 1)
$cell.Range.Find.Execute("11÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "61÷3=", 1)

$cell = $tbl.Cell(5, 2)
$cell.Range.Find.Execute("80÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "20÷4=", 1)

$cell = $tbl.Cell(5, 3)
$cell.Range.Find.Execute("38÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "13÷6=", 1)

$cell = $tbl.Cell(5, 4)
$cell.Range.Find.Execute("77÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷5=", 1)

$cell = $tbl.Cell(5, 5)
$cell.Range.Find.Execute("49÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷3=", 1)

$cell = $tbl.Cell(9, 1)
$cell.Range.Find.Execute("40÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "85÷8=", 1)

$cell = $tbl.Cell(9, 2)
$cell.Range.Find.Execute("14÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷4=", 1)

$cell = $tbl.Cell(9, 3)
$cell.Range.Find.Execute("62÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "24÷4=", 1)

$cell = $tbl.Cell(9, 4)
$cell.Range.Find.Execute("56÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "49÷8=", 1)

$cell = $tbl.Cell(9, 5)
$cell.Range.Find.Execute("79÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "85÷9=", 1)

$cell = $tbl.Cell(13, 1)
$cell.Range.Find.Execute("56÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷9=", 1)

$cell = $tbl.Cell(13, 2)
$cell.Range.Find.Execute("45÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "11÷5=", 1)

$cell = $tbl.Cell(13, 3)
$cell.Range.Find.Execute("18÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "39÷4=", 1)

$cell = $tbl.Cell(13, 4)
$cell.Range.Find.Execute("91÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "23÷7=", 1)

$cell = $tbl.Cell(13, 5)
$cell.Range.Find.Execute("96÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷6=", 1)

$cell = $tbl.Cell(17, 1)
$cell.Range.Find.Execute("59÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "38÷6=", 1)

$cell = $tbl.Cell(17, 2)
$cell.Range.Find.Execute("29÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "27÷7=", 1)

$cell = $tbl.Cell(17, 3)
$cell.Range.Find.Execute("55÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "33÷3=", 1)

$cell = $tbl.Cell(17, 4)
$cell.Range.Find.Execute("38÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷3=", 1)

$cell = $tbl.Cell(17, 5)
$cell.Range.Find.Execute("38÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "60÷3=", 1)
